# Adds 9 new data rows (22-30) to the master-reg_center_machine_h sheet,
# mirroring the existing row pattern (regcntr_id, machine_id, lang_code,
# is_active, cr_by, cr_dtimes, eff_dtimes), updates the selection to the
# newly added machine_id block, and sets the page to portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id (col A) and machine_id (col B) values for the new rows 22-30
$regCenterIds = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)
$machineIds   = @(10021, 10022, 10023, 10024, 10025, 10026, 10027, 10028, 10029)

$startRow = 22
for ($i = 0; $i -lt $regCenterIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $regCenterIds[$i]
    $ws.Cells.Item($row, 2).Value = $machineIds[$i]
}

# Columns shared across every new row: lang_code, is_active, cr_by, cr_dtimes, eff_dtimes
$ws.Range("C22:C30").Value = "eng"
$ws.Range("D22:D30").Value = $true
$ws.Range("E22:E30").Value = "superadmin"
$ws.Range("F22:G30").Value = "now()"

# Match the workbook's saved selection state (machine_id column of the new rows)
$ws.Range("B22:B30").Select()

# Print setup: portrait orientation
$ws.PageSetup.Orientation = 1
